$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 238.6
$ws.Range("I18").Value = 238.6
$ws.Range("K18").Value = 238.6
$ws.Range("M18").Value = 45.40000000000001

$ws.Range("H76").Value = 3401.9048
$ws.Range("I76").Value = 2966.4707
$ws.Range("K76").Value = 2966.4707
$ws.Range("M76").Value = -2651.4707

$ws.Range("H79").Value = 3401.9048
$ws.Range("I79").Value = 2966.4707
$ws.Range("K79").Value = 2966.4707
$ws.Range("M79").Value = -1874.4707

$ws.Range("H116").Value = 3323.3
$ws.Range("J116").Value = 4356.154
$ws.Range("L116").Value = 4356.154
$ws.Range("N116").Value = -11240.154

$ws.Range("H129").Value = 1030.8718
$ws.Range("J129").Value = 1045.44
$ws.Range("L129").Value = 3136.32
$ws.Range("N129").Value = -13136.32

$ws.Range("H130").Value = 14504.4
$ws.Range("I130").Value = 2700
$ws.Range("J130").Value = 14996.25
$ws.Range("K130").Value = 2700
$ws.Range("L130").Value = 14996.25
$ws.Range("M130").Value = 2320
$ws.Range("N130").Value = -25036.25

$ws.Range("H137").Value = 1342.1482
$ws.Range("I137").Value = 1571.8182
$ws.Range("J137").Value = 1184.25
$ws.Range("K137").Value = 4715.4546
$ws.Range("L137").Value = 3552.75
$ws.Range("M137").Value = -2165.4546
$ws.Range("N137").Value = -8652.75

$ws.Range("H138").Value = 3975.7092
$ws.Range("I138").Value = 2052.6
$ws.Range("J138").Value = 6086.439
$ws.Range("K138").Value = 6157.799999999999
$ws.Range("L138").Value = 18259.317
$ws.Range("M138").Value = -1017.799999999999
$ws.Range("N138").Value = -28539.317

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 29920
$ws.Range("J134").Value = 31573.334
$ws.Range("L134").Value = 31573.334
$ws.Range("N134").Value = -41713.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 1722.5111
$ws.Range("I134").Value = 1298.9688
$ws.Range("J134").Value = 2765.077
$ws.Range("K134").Value = 3896.9064
$ws.Range("L134").Value = 8295.231
$ws.Range("M134").Value = -1361.9064
$ws.Range("N134").Value = -13365.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3969.052
$ws.Range("I31").Value = 2082.0735
$ws.Range("K31").Value = 2082.0735
$ws.Range("M31").Value = -1787.0735

$ws.Range("H34").Value = 3969.052
$ws.Range("I34").Value = 2082.0735
$ws.Range("K34").Value = 2082.0735
$ws.Range("M34").Value = -1880.0735

$ws.Range("H38").Value = 12719
$ws.Range("I38").Value = 12719
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 12719
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -12342
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 12719
$ws.Range("I46").Value = 12719
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 12719
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -12508
$ws.Range("N46").ClearContents()

$ws.Range("H50").Value = 30156
$ws.Range("J50").Value = 30156
$ws.Range("L50").Value = 30156
$ws.Range("N50").Value = -31406

$ws.Range("H51").Value = 30795.096
$ws.Range("J51").Value = 31834.85
$ws.Range("L51").Value = 31834.85
$ws.Range("N51").Value = -33306.85

$ws.Range("H58").Value = 3283.2727
$ws.Range("I58").Value = 3159.1428
$ws.Range("J58").Value = 3341.2
$ws.Range("K58").Value = 3159.1428
$ws.Range("L58").Value = 3341.2
$ws.Range("M58").Value = -2956.1428
$ws.Range("N58").Value = -3747.2

$ws.Range("H59").Value = 30960
$ws.Range("J59").Value = 33700
$ws.Range("L59").Value = 33700
$ws.Range("N59").Value = -35990

$ws.Range("H61").Value = 30795.096
$ws.Range("J61").Value = 31834.85
$ws.Range("L61").Value = 31834.85
$ws.Range("N61").Value = -32530.85

$ws.Range("H74").Value = 39000
$ws.Range("J74").Value = 39000
$ws.Range("L74").Value = 39000
$ws.Range("N74").Value = -40748

$ws.Range("H77").Value = 39000
$ws.Range("J77").Value = 39000
$ws.Range("L77").Value = 117000
$ws.Range("N77").Value = -125736

$ws.Range("H136").Value = 3283.2727
$ws.Range("I136").Value = 3159.1428
$ws.Range("J136").Value = 3341.2
$ws.Range("K136").Value = 9477.428400000001
$ws.Range("L136").Value = 10023.6
$ws.Range("M136").Value = -6927.428400000001
$ws.Range("N136").Value = -15123.6

$ws.Range("H140").Value = 74209.55499999999
$ws.Range("J140").Value = 74209.55499999999
$ws.Range("L140").Value = 74209.55499999999
$ws.Range("N140").Value = -84569.55499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 825
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 4500
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -6996

$ws.Range("H94").Value = 6254.5
$ws.Range("I94").Value = 800
$ws.Range("J94").Value = 7345.4
$ws.Range("K94").Value = 2400
$ws.Range("L94").Value = 22036.2
$ws.Range("M94").Value = -1724
$ws.Range("N94").Value = -23388.2

$ws.Range("H131").Value = 775.78
$ws.Range("I131").Value = 450.11765
$ws.Range("J131").Value = 842.48193
$ws.Range("K131").Value = 1350.35295
$ws.Range("L131").Value = 2527.44579
$ws.Range("M131").Value = 3689.64705
$ws.Range("N131").Value = -12607.44579

$ws.Range("H137").Value = 7833.5405
$ws.Range("I137").Value = 2978.4285
$ws.Range("J137").Value = 8966.4
$ws.Range("K137").Value = 8935.2855
$ws.Range("L137").Value = 26899.2
$ws.Range("M137").Value = -3835.2855
$ws.Range("N137").Value = -37099.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2400.182
$ws.Range("I126").Value = 2333
$ws.Range("J126").Value = 2517.75
$ws.Range("K126").Value = 6999
$ws.Range("L126").Value = 7553.25
$ws.Range("M126").Value = -4529
$ws.Range("N126").Value = -12493.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1577.5714
$ws.Range("I46").Value = 1677.2
$ws.Range("J46").Value = 1522.2222
$ws.Range("K46").Value = 1677.2
$ws.Range("L46").Value = 1522.2222
$ws.Range("M46").Value = -1489.2
$ws.Range("N46").Value = -1898.2222

$ws.Range("H55").Value = 262.23077
$ws.Range("I55").Value = 141.66667
$ws.Range("J55").Value = 365.57144
$ws.Range("K55").Value = 141.66667
$ws.Range("L55").Value = 365.57144
$ws.Range("M55").Value = 31.33332999999999
$ws.Range("N55").Value = -711.5714399999999

$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19900

$ws.Range("H127").Value = 57194
$ws.Range("J127").Value = 57194
$ws.Range("L127").Value = 57194
$ws.Range("N127").Value = -67114

$ws.Range("H132").Value = 10446.318
$ws.Range("I132").Value = 3839.1875
$ws.Range("J132").Value = 28065.334
$ws.Range("K132").Value = 11517.5625
$ws.Range("L132").Value = 84196.00199999999
$ws.Range("M132").Value = -8987.5625
$ws.Range("N132").Value = -89256.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 55070.4
$ws.Range("J137").Value = 55070.4
$ws.Range("L137").Value = 55070.4
$ws.Range("N137").Value = -65270.4
